# Update cached market-price / profit figures on the Leve-profit sheets.
# Values sourced from the latest scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 28
$ws.Cells.Item(28, 8).Value = 347.5  # H28: 339.7 -> 347.5
$ws.Cells.Item(28, 9).Value = 245.72728  # I28: 232.66667 -> 245.72728
$ws.Cells.Item(28, 10).Value = 507.42856  # J28: 500.25 -> 507.42856
$ws.Cells.Item(28, 11).Value = 245.72728  # K28: 232.66667 -> 245.72728
$ws.Cells.Item(28, 12).Value = 507.42856  # L28: 500.25 -> 507.42856
$ws.Cells.Item(28, 13).Value = 239.27272  # M28: 252.33333 -> 239.27272
$ws.Cells.Item(28, 14).Value = -1477.42856  # N28: -1470.25 -> -1477.42856

# Row 113
$ws.Cells.Item(113, 8).Value = 10418742  # H113: 15627200 -> 10418742
$ws.Cells.Item(113, 9).Value = 2234  # I113: 2475 -> 2234
$ws.Cells.Item(113, 10).Value = 20835250  # J113: 31251924 -> 20835250
$ws.Cells.Item(113, 11).Value = 2234  # K113: 2475 -> 2234
$ws.Cells.Item(113, 12).Value = 20835250  # L113: 31251924 -> 20835250
$ws.Cells.Item(113, 13).Value = 1020  # M113: 779 -> 1020
$ws.Cells.Item(113, 14).Value = -20841758  # N113: -31258432 -> -20841758

# Row 135
$ws.Cells.Item(135, 8).Value = 4109  # H135: 1344.7542 -> 4109
$ws.Cells.Item(135, 9).Value = 3764.375  # I135: 1356.4576 -> 3764.375
$ws.Cells.Item(135, 10).Value = 8244.5  # J135: 999.5 -> 8244.5
$ws.Cells.Item(135, 11).Value = 33879.375  # K135: 12208.1184 -> 33879.375
$ws.Cells.Item(135, 12).Value = 74200.5  # L135: 8995.5 -> 74200.5
$ws.Cells.Item(135, 13).Value = -31344.375  # M135: -9673.118399999999 -> -31344.375
$ws.Cells.Item(135, 14).Value = -79270.5  # N135: -14065.5 -> -79270.5

# Row 137
$ws.Cells.Item(137, 8).Value = 1926.0834  # H137: 2226.158 -> 1926.0834
$ws.Cells.Item(137, 9).Value = 1886.3  # I137: 2253.1333 -> 1886.3
$ws.Cells.Item(137, 11).Value = 5658.9  # K137: 6759.3999 -> 5658.9
$ws.Cells.Item(137, 13).Value = -3108.9  # M137: -4209.3999 -> -3108.9

$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Cells.Item(32, 8).Value = 14409.293  # H32: 14038.813 -> 14409.293
$ws.Cells.Item(32, 9).Value = 10225.71  # I32: 9966.212 -> 10225.71
$ws.Cells.Item(32, 10).Value = 27378.4  # J32: 27478.4 -> 27378.4
$ws.Cells.Item(32, 11).Value = 10225.71  # K32: 9966.212 -> 10225.71
$ws.Cells.Item(32, 12).Value = 27378.4  # L32: 27478.4 -> 27378.4
$ws.Cells.Item(32, 13).Value = -9938.709999999999  # M32: -9679.212 -> -9938.709999999999
$ws.Cells.Item(32, 14).Value = -27952.4  # N32: -28052.4 -> -27952.4

# Row 45
$ws.Cells.Item(45, 8).Value = 211276.4  # H45: 117825.22 -> 211276.4
$ws.Cells.Item(45, 9).Value = 350544  # I45: 175612.83 -> 350544
$ws.Cells.Item(45, 10).Value = 2375  # J45: 2250 -> 2375
$ws.Cells.Item(45, 11).Value = 350544  # K45: 175612.83 -> 350544
$ws.Cells.Item(45, 12).Value = 2375  # L45: 2250 -> 2375
$ws.Cells.Item(45, 13).Value = -350167  # M45: -175235.83 -> -350167
$ws.Cells.Item(45, 14).Value = -3129  # N45: -3004 -> -3129

# Row 74
$ws.Cells.Item(74, 8).Value = 1126.597  # H74: 1090.5071 -> 1126.597
$ws.Cells.Item(74, 9).Value = 796.9048  # I74: 769.86957 -> 796.9048
$ws.Cells.Item(74, 11).Value = 796.9048  # K74: 769.86957 -> 796.9048
$ws.Cells.Item(74, 13).Value = 77.09519999999998  # M74: 104.13043 -> 77.09519999999998

# Row 77
$ws.Cells.Item(77, 8).Value = 1126.597  # H77: 1090.5071 -> 1126.597
$ws.Cells.Item(77, 9).Value = 796.9048  # I77: 769.86957 -> 796.9048
$ws.Cells.Item(77, 11).Value = 3984.524  # K77: 3849.34785 -> 3984.524
$ws.Cells.Item(77, 13).Value = 383.4759999999997  # M77: 518.6521500000003 -> 383.4759999999997

# Row 102
$ws.Cells.Item(102, 8).Value = 3089370.5  # H102: 3706963.8 -> 3089370.5
$ws.Cells.Item(102, 9).Value = 3705344.5  # I102: 4631329.5 -> 3705344.5
$ws.Cells.Item(102, 11).Value = 3705344.5  # K102: 4631329.5 -> 3705344.5
$ws.Cells.Item(102, 13).Value = -3703722.5  # M102: -4629707.5 -> -3703722.5

$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Cells.Item(86, 8).Value = 1798.7273  # H86: 1847.2 -> 1798.7273
$ws.Cells.Item(86, 9).Value = 1816.3334  # I86: 1745.3334 -> 1816.3334
$ws.Cells.Item(86, 10).Value = 1777.6  # J86: 2000 -> 1777.6
$ws.Cells.Item(86, 11).Value = 1816.3334  # K86: 1745.3334 -> 1816.3334
$ws.Cells.Item(86, 12).Value = 1777.6  # L86: 2000 -> 1777.6
$ws.Cells.Item(86, 13).Value = -693.3334  # M86: -622.3334 -> -693.3334
$ws.Cells.Item(86, 14).Value = -4023.6  # N86: -4246 -> -4023.6

# Row 89
$ws.Cells.Item(89, 8).Value = 1798.7273  # H89: 1847.2 -> 1798.7273
$ws.Cells.Item(89, 9).Value = 1816.3334  # I89: 1745.3334 -> 1816.3334
$ws.Cells.Item(89, 10).Value = 1777.6  # J89: 2000 -> 1777.6
$ws.Cells.Item(89, 11).Value = 9081.666999999999  # K89: 8726.666999999999 -> 9081.666999999999
$ws.Cells.Item(89, 12).Value = 8888  # L89: 10000 -> 8888
$ws.Cells.Item(89, 13).Value = -3465.666999999999  # M89: -3110.666999999999 -> -3465.666999999999
$ws.Cells.Item(89, 14).Value = -20120  # N89: -21232 -> -20120

# Row 107
$ws.Cells.Item(107, 8).Value = 500970.34  # H107: 429603.16 -> 500970.34
$ws.Cells.Item(107, 9).Value = 500970.34  # I107: 601064.4 -> 500970.34
$ws.Cells.Item(107, 10).Value = 0  # J107: 950 -> 0
$ws.Cells.Item(107, 11).Value = 500970.34  # K107: 601064.4 -> 500970.34
$ws.Cells.Item(107, 12).Value = 0  # L107: 950 -> 0
$ws.Cells.Item(107, 13).Value = -499050.34  # M107: -599144.4 -> -499050.34
$ws.Cells.Item(107, 14).ClearContents()  # N107: -4790 -> (removed)

# Row 134
$ws.Cells.Item(134, 8).Value = 22877.9  # H134: 25832.887 -> 22877.9
$ws.Cells.Item(134, 9).Value = 2877.3408  # I134: 3187.973 -> 2877.3408
$ws.Cells.Item(134, 10).Value = 169548.67  # J134: 145527.42 -> 169548.67
$ws.Cells.Item(134, 11).Value = 8632.0224  # K134: 9563.919 -> 8632.0224
$ws.Cells.Item(134, 12).Value = 508646.01  # L134: 436582.26 -> 508646.01
$ws.Cells.Item(134, 13).Value = -6097.0224  # M134: -7028.919 -> -6097.0224
$ws.Cells.Item(134, 14).Value = -513716.01  # N134: -441652.26 -> -513716.01

$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Cells.Item(16, 8).Value = 4527084.5  # H16: 4275535 -> 4527084.5
$ws.Cells.Item(16, 9).Value = 8548867  # I16: 5918845.5 -> 8548867
$ws.Cells.Item(16, 10).Value = 2579.375  # J16: 2927.8 -> 2579.375
$ws.Cells.Item(16, 11).Value = 8548867  # K16: 5918845.5 -> 8548867
$ws.Cells.Item(16, 12).Value = 2579.375  # L16: 2927.8 -> 2579.375
$ws.Cells.Item(16, 13).Value = -8548580  # M16: -5918558.5 -> -8548580
$ws.Cells.Item(16, 14).Value = -3153.375  # N16: -3501.8 -> -3153.375

# Row 31
$ws.Cells.Item(31, 8).Value = 6973.4717  # H31: 6639.9243 -> 6973.4717
$ws.Cells.Item(31, 9).Value = 1678.2  # I31: 1478.5834 -> 1678.2
$ws.Cells.Item(31, 10).Value = 10182.728  # J31: 10911.379 -> 10182.728
$ws.Cells.Item(31, 11).Value = 1678.2  # K31: 1478.5834 -> 1678.2
$ws.Cells.Item(31, 12).Value = 10182.728  # L31: 10911.379 -> 10182.728
$ws.Cells.Item(31, 13).Value = -1383.2  # M31: -1183.5834 -> -1383.2
$ws.Cells.Item(31, 14).Value = -10772.728  # N31: -11501.379 -> -10772.728

# Row 34
$ws.Cells.Item(34, 8).Value = 6973.4717  # H34: 6639.9243 -> 6973.4717
$ws.Cells.Item(34, 9).Value = 1678.2  # I34: 1478.5834 -> 1678.2
$ws.Cells.Item(34, 10).Value = 10182.728  # J34: 10911.379 -> 10182.728
$ws.Cells.Item(34, 11).Value = 1678.2  # K34: 1478.5834 -> 1678.2
$ws.Cells.Item(34, 12).Value = 10182.728  # L34: 10911.379 -> 10182.728
$ws.Cells.Item(34, 13).Value = -1476.2  # M34: -1276.5834 -> -1476.2
$ws.Cells.Item(34, 14).Value = -10586.728  # N34: -11315.379 -> -10586.728

# Row 105
$ws.Cells.Item(105, 8).Value = 33336834  # H105: 17546004 -> 33336834
$ws.Cells.Item(105, 9).Value = 66671530  # I105: 25643376 -> 66671530
$ws.Cells.Item(105, 10).Value = 2140  # J105: 1700 -> 2140
$ws.Cells.Item(105, 11).Value = 66671530  # K105: 25643376 -> 66671530
$ws.Cells.Item(105, 12).Value = 2140  # L105: 1700 -> 2140
$ws.Cells.Item(105, 13).Value = -66669783  # M105: -25641629 -> -66669783
$ws.Cells.Item(105, 14).Value = -5634  # N105: -5194 -> -5634

# Row 113
$ws.Cells.Item(113, 8).Value = 4527084.5  # H113: 4275535 -> 4527084.5
$ws.Cells.Item(113, 9).Value = 8548867  # I113: 5918845.5 -> 8548867
$ws.Cells.Item(113, 10).Value = 2579.375  # J113: 2927.8 -> 2579.375
$ws.Cells.Item(113, 11).Value = 8548867  # K113: 5918845.5 -> 8548867
$ws.Cells.Item(113, 12).Value = 2579.375  # L113: 2927.8 -> 2579.375
$ws.Cells.Item(113, 13).Value = -8546697  # M113: -5916675.5 -> -8546697
$ws.Cells.Item(113, 14).Value = -6919.375  # N113: -7267.8 -> -6919.375

$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Cells.Item(3, 8).Value = 8303.808000000001  # H3: 10614.223 -> 8303.808000000001
$ws.Cells.Item(3, 9).Value = 11791.667  # I3: 21165 -> 11791.667
$ws.Cells.Item(3, 10).Value = 5314.2144  # J3: 5338.8335 -> 5314.2144
$ws.Cells.Item(3, 11).Value = 35375.001  # K3: 63495 -> 35375.001
$ws.Cells.Item(3, 12).Value = 15942.6432  # L3: 16016.5005 -> 15942.6432
$ws.Cells.Item(3, 13).Value = -35263.001  # M3: -63383 -> -35263.001
$ws.Cells.Item(3, 14).Value = -16166.6432  # N3: -16240.5005 -> -16166.6432

# Row 131
$ws.Cells.Item(131, 8).Value = 1755487.6  # H131: 1961867 -> 1755487.6
$ws.Cells.Item(131, 9).Value = 4762439.5  # I131: 5000556 -> 4762439.5
$ws.Cells.Item(131, 10).Value = 1432.2778  # J131: 1422.5807 -> 1432.2778
$ws.Cells.Item(131, 11).Value = 14287318.5  # K131: 15001668 -> 14287318.5
$ws.Cells.Item(131, 12).Value = 4296.8334  # L131: 4267.742099999999 -> 4296.8334
$ws.Cells.Item(131, 13).Value = -14282278.5  # M131: -14996628 -> -14282278.5
$ws.Cells.Item(131, 14).Value = -14376.8334  # N131: -14347.7421 -> -14376.8334

# Row 132
$ws.Cells.Item(132, 8).Value = 2772.111  # H132: 4631979.5 -> 2772.111
$ws.Cells.Item(132, 9).Value = 1399.125  # I132: 1369.3 -> 1399.125
$ws.Cells.Item(132, 10).Value = 3870.5  # J132: 7939558.5 -> 3870.5
$ws.Cells.Item(132, 11).Value = 12592.125  # K132: 12323.7 -> 12592.125
$ws.Cells.Item(132, 12).Value = 34834.5  # L132: 71456026.5 -> 34834.5
$ws.Cells.Item(132, 13).Value = -10062.125  # M132: -9793.699999999999 -> -10062.125
$ws.Cells.Item(132, 14).Value = -39894.5  # N132: -71461086.5 -> -39894.5

# Row 133
$ws.Cells.Item(133, 8).Value = 39696.1  # H133: 63288.832 -> 39696.1
$ws.Cells.Item(133, 9).Value = 116818.89  # I133: 131943.75 -> 116818.89
$ws.Cells.Item(133, 10).Value = 6643.476  # J133: 8364.9 -> 6643.476
$ws.Cells.Item(133, 11).Value = 350456.67  # K133: 395831.25 -> 350456.67
$ws.Cells.Item(133, 12).Value = 19930.428  # L133: 25094.7 -> 19930.428
$ws.Cells.Item(133, 13).Value = -345396.67  # M133: -390771.25 -> -345396.67
$ws.Cells.Item(133, 14).Value = -30050.428  # N133: -35214.7 -> -30050.428

# Row 137
$ws.Cells.Item(137, 8).Value = 21153.316  # H137: 17676.791 -> 21153.316
$ws.Cells.Item(137, 9).Value = 9777.143  # I137: 9854.286 -> 9777.143
$ws.Cells.Item(137, 10).Value = 53006.6  # J137: 28628.3 -> 53006.6
$ws.Cells.Item(137, 11).Value = 29331.429  # K137: 29562.858 -> 29331.429
$ws.Cells.Item(137, 12).Value = 159019.8  # L137: 85884.89999999999 -> 159019.8
$ws.Cells.Item(137, 13).Value = -24231.429  # M137: -24462.858 -> -24231.429
$ws.Cells.Item(137, 14).Value = -169219.8  # N137: -96084.89999999999 -> -169219.8

# Row 139
$ws.Cells.Item(139, 8).Value = 4792.394  # H139: 5050.5938 -> 4792.394
$ws.Cells.Item(139, 9).Value = 5957.45  # I139: 6645.5 -> 5957.45
$ws.Cells.Item(139, 11).Value = 17872.35  # K139: 19936.5 -> 17872.35
$ws.Cells.Item(139, 13).Value = -12732.35  # M139: -14796.5 -> -12732.35

# Row 140
$ws.Cells.Item(140, 8).Value = 4743.9443  # H140: 5496.4 -> 4743.9443
$ws.Cells.Item(140, 9).Value = 4743.9443  # I140: 5496.4 -> 4743.9443
$ws.Cells.Item(140, 11).Value = 14231.8329  # K140: 16489.2 -> 14231.8329
$ws.Cells.Item(140, 13).Value = -9051.832900000001  # M140: -11309.2 -> -9051.832900000001

$ws = $wb.Worksheets.Item("GSM")

# Row 122
$ws.Cells.Item(122, 8).Value = 32872316  # H122: 36888908 -> 32872316
$ws.Cells.Item(122, 9).Value = 36719244  # I122: 38030772 -> 36719244
$ws.Cells.Item(122, 10).Value = 22730412  # J122: 33336436 -> 22730412
$ws.Cells.Item(122, 11).Value = 110157732  # K122: 114092316 -> 110157732
$ws.Cells.Item(122, 12).Value = 68191236  # L122: 100009308 -> 68191236
$ws.Cells.Item(122, 13).Value = -110155282  # M122: -114089866 -> -110155282
$ws.Cells.Item(122, 14).Value = -68196136  # N122: -100014208 -> -68196136

# Row 132
$ws.Cells.Item(132, 8).Value = 14055.177  # H132: 16680989 -> 14055.177
$ws.Cells.Item(132, 9).Value = 9795  # I132: 20844786 -> 9795
$ws.Cells.Item(132, 10).Value = 46006.5  # J132: 25799.5 -> 46006.5
$ws.Cells.Item(132, 11).Value = 29385  # K132: 62534358 -> 29385
$ws.Cells.Item(132, 12).Value = 138019.5  # L132: 77398.5 -> 138019.5
$ws.Cells.Item(132, 13).Value = -26855  # M132: -62531828 -> -26855
$ws.Cells.Item(132, 14).Value = -143079.5  # N132: -82458.5 -> -143079.5

$ws = $wb.Worksheets.Item("LTW")

# Row 127
$ws.Cells.Item(127, 8).Value = 37500  # H127: 0 -> 37500
$ws.Cells.Item(127, 10).Value = 37500  # J127: 0 -> 37500
$ws.Cells.Item(127, 12).Value = 37500  # L127: 0 -> 37500
$ws.Cells.Item(127, 14).Value = -47420  # N127: (new) -> -47420

# Row 132
$ws.Cells.Item(132, 8).Value = 9093433  # H132: 11630684 -> 9093433
$ws.Cells.Item(132, 9).Value = 10002396  # I132: 12197839 -> 10002396
$ws.Cells.Item(132, 10).Value = 3799  # J132: 4002.5 -> 3799
$ws.Cells.Item(132, 11).Value = 30007188  # K132: 36593517 -> 30007188
$ws.Cells.Item(132, 12).Value = 11397  # L132: 12007.5 -> 11397
$ws.Cells.Item(132, 13).Value = -30004658  # M132: -36590987 -> -30004658
$ws.Cells.Item(132, 14).Value = -16457  # N132: -17067.5 -> -16457

$ws = $wb.Worksheets.Item("WVR")

# Row 136
$ws.Cells.Item(136, 8).Value = 4868.222  # H136: 5559190 -> 4868.222
$ws.Cells.Item(136, 9).Value = 5698  # I136: 4458.533 -> 5698
$ws.Cells.Item(136, 10).Value = 3564.2856  # J136: 11113921 -> 3564.2856
$ws.Cells.Item(136, 11).Value = 17094  # K136: 13375.599 -> 17094
$ws.Cells.Item(136, 12).Value = 10692.8568  # L136: 33341763 -> 10692.8568
$ws.Cells.Item(136, 13).Value = -14544  # M136: -10825.599 -> -14544
$ws.Cells.Item(136, 14).Value = -15792.8568  # N136: -33346863 -> -15792.8568
